$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the case data paths: testdata/mysqlcases -> testdata/btreemysqlcases
# for rows M3, M4, M5 (row M2 already uses btreemysqlcases)
$ws.Range("M3").Value = "src/test/resources/io.dingodb.test/testdata/btreemysqlcases/prepareStatement/expectedresult/btree_ps_dml_002.csv"
$ws.Range("M4").Value = "src/test/resources/io.dingodb.test/testdata/btreemysqlcases/prepareStatement/expectedresult/btree_ps_dml_003.csv"
$ws.Range("M5").Value = "src/test/resources/io.dingodb.test/testdata/btreemysqlcases/prepareStatement/expectedresult/btree_ps_dml_004.csv"

# Update the selected cell/active selection on Sheet1 to L13
$ws.Activate()
$ws.Range("L13").Select()
